$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look numeric to Excel auto-detection; force them
# to stay as plain text (matching the original inline-string cell type) by
# pre-applying a text number format before writing the value.
$textCells = @(
    'D5',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D15',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D25',
    'D26',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = @(
    @('D2', '27.935.73'),
    @('E2', '  -0.28%  '),
    @('D3', '1.856.80'),
    @('E3', '  -1.46%  '),
    @('E4', '  +0.36%  '),
    @('D5', '311.35'),
    @('E5', '  -0.46%  '),
    @('E6', '  +0.25%  '),
    @('D7', '0.5112'),
    @('E7', '  +2.33%  '),
    @('D8', '0.3797'),
    @('E8', '  -1.56%  '),
    @('D9', '0.08264'),
    @('E9', '  -10.03%  '),
    @('D10', '41.58'),
    @('E10', '  -0.21%  '),
    @('D11', '1.105'),
    @('E11', '  -1.52%  '),
    @('D12', '6.176'),
    @('E12', '  -2.46%  '),
    @('D13', '1.857.41'),
    @('E13', '  -1.34%  '),
    @('D14', '20.38'),
    @('E14', '  -1.77%  '),
    @('D15', '7.170'),
    @('E15', '  -1.61%  '),
    @('E16', '  +0.27%  '),
    @('D17', '0.00001092'),
    @('E17', '  -1.19%  '),
    @('D18', '90.09'),
    @('E18', '  -1.44%  '),
    @('D19', '0.06606'),
    @('E19', '  -0.45%  '),
    @('D20', '17.64'),
    @('E20', '  -1.91%  '),
    @('D21', '1.002'),
    @('D22', '5.994'),
    @('E22', '  -3.11%  '),
    @('D23', '27.972.98'),
    @('E23', '  -0.26%  '),
    @('D24', '10.98'),
    @('E24', '  -3.66%  '),
    @('D25', '2.225'),
    @('E25', '  -3.33%  '),
    @('D26', '2.568'),
    @('E26', '  +0.87%  '),
    @('D27', '2.068.94'),
    @('E27', '  -1.51%  '),
    @('D28', '156.83'),
    @('E28', '  -0.48%  '),
    @('D29', '20.36'),
    @('E29', '  -1.96%  '),
    @('D30', '124.30'),
    @('E30', '  -1.85%  '),
    @('D31', '0.1058'),
    @('E31', '  +0.15%  '),
    @('D32', '1.035'),
    @('E32', '  -3.05%  '),
    @('D33', '5.574'),
    @('E33', '  -0.32%  '),
    @('D34', '3.602'),
    @('E34', '  +0.59%  '),
    @('D35', '9.491'),
    @('E35', '  +1.37%  '),
    @('D36', '0.06497'),
    @('E36', '  -1.27%  '),
    @('D37', '0.02400'),
    @('E37', '  -0.12%  '),
    @('D38', '0.2150'),
    @('E38', '  -1.84%  '),
    @('D39', '1.201'),
    @('E39', '  -0.63%  '),
    @('D40', '0.6397'),
    @('E40', '  -0.23%  '),
    @('D41', '1.224'),
    @('E41', '  -5.08%  '),
    @('D42', '11.21'),
    @('E42', '  -2.95%  '),
    @('D43', '4.849'),
    @('E43', '  -1.69%  '),
    @('D44', '0.6059'),
    @('E44', '  +0.28%  '),
    @('D45', '13.01'),
    @('E45', '  -2.13%  '),
    @('D46', '1.280'),
    @('E46', '  -0.97%  '),
    @('D47', '3.664'),
    @('E47', '  -0.14%  '),
    @('E48', '  -1.26%  '),
    @('E49', '  -0.64%  '),
    @('D50', '120.52'),
    @('E50', '  -0.59%  '),
    @('D51', '79.73'),
    @('E51', '  +1.26%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
